# "get patient age and sex"
#
# The patient row (row 12) previously pointed its age/sex lookup at the
# old sagittal-patient-data workbook; point it at the new consolidated
# patient-information workbook instead. The "manual_landmarks" column
# (L12) keeps its literal value ("_") throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPath = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\patient_information.xlsx"

$ws.Range("G12").Value = $newPath
$ws.Range("H12").Value = $newPath

# Update the current selection to reflect where the author ended up
# working (J13) after making the edit.
$ws.Range("J13").Select()
